# Auto-generated Excel COM-interop script applying the Zalera_Profits diff
# Updates the Leve profit-calculation columns (H,I,J,K,L,M,N) for specific
# rows across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching
# the scheduled-runner data refresh described by the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H40").Value = 3294
$ws.Range("I40").Value = 3930.75
$ws.Range("K40").Value = 3930.75
$ws.Range("M40").Value = -3755.75

$ws.Range("H76").Value = 7150699.5
$ws.Range("I76").Value = 12506230
$ws.Range("K76").Value = 12506230
$ws.Range("M76").Value = -12505915

$ws.Range("H79").Value = 7150699.5
$ws.Range("I79").Value = 12506230
$ws.Range("K79").Value = 12506230
$ws.Range("M79").Value = -12505138

$ws.Range("H80").Value = 1425
$ws.Range("I80").Value = 1201.1111
$ws.Range("K80").Value = 3603.3333
$ws.Range("M80").Value = -2605.3333

$ws.Range("H83").Value = 1425
$ws.Range("I83").Value = 1201.1111
$ws.Range("K83").Value = 10809.9999
$ws.Range("M83").Value = -5817.999900000001

$ws.Range("H92").Value = 1561.4166
$ws.Range("I92").Value = 1561.4166
$ws.Range("K92").Value = 1561.4166
$ws.Range("M92").Value = -313.4166

$ws.Range("H98").Value = 5214.409
$ws.Range("I98").Value = 4703.1113
$ws.Range("K98").Value = 4703.1113
$ws.Range("M98").Value = -3205.1113

$ws.Range("H116").Value = 3995
$ws.Range("J116").Value = 3995
$ws.Range("L116").Value = 3995
$ws.Range("N116").Value = -10879

$ws.Range("H122").Value = 5214.409
$ws.Range("I122").Value = 4703.1113
$ws.Range("K122").Value = 14109.3339
$ws.Range("M122").Value = -11659.3339

$ws.Range("H125").Value = 2862.375
$ws.Range("I125").Value = 2029.8
$ws.Range("J125").Value = 4250
$ws.Range("K125").Value = 18268.2
$ws.Range("L125").Value = 38250
$ws.Range("M125").Value = -15808.2
$ws.Range("N125").Value = -43170

$ws.Range("H137").Value = 22731826
$ws.Range("I137").Value = 250000000
$ws.Range("J137").Value = 5009.6
$ws.Range("K137").Value = 750000000
$ws.Range("L137").Value = 15028.8
$ws.Range("M137").Value = -749997450
$ws.Range("N137").Value = -20128.8


$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 41781.57
$ws.Range("I32").Value = 44787.848
$ws.Range("K32").Value = 44787.848
$ws.Range("M32").Value = -44500.848

$ws.Range("H44").Value = 76249.5
$ws.Range("J44").Value = 76249.5
$ws.Range("L44").Value = 76249.5
$ws.Range("N44").Value = -77225.5

$ws.Range("H55").Value = 8500
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H61").Value = 7392.909
$ws.Range("J61").Value = 10663.2
$ws.Range("L61").Value = 10663.2
$ws.Range("N61").Value = -11087.2

$ws.Range("H132").Value = 4338.275
$ws.Range("I132").Value = 3097.7273
$ws.Range("K132").Value = 9293.1819
$ws.Range("M132").Value = -6763.1819

$ws.Range("H136").Value = 7392.909
$ws.Range("J136").Value = 10663.2
$ws.Range("L136").Value = 31989.6
$ws.Range("N136").Value = -37089.60000000001


$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H6").Value = 45000
$ws.Range("J6").Value = 45000
$ws.Range("L6").Value = 45000
$ws.Range("N6").Value = -45226

$ws.Range("H20").Value = 2931.2222
$ws.Range("I20").Value = 1999.6666
$ws.Range("K20").Value = 1999.6666
$ws.Range("M20").Value = -1752.6666

$ws.Range("H86").Value = 644856.7
$ws.Range("I86").Value = 2999
$ws.Range("K86").Value = 2999
$ws.Range("M86").Value = -1876

$ws.Range("H89").Value = 644856.7
$ws.Range("I89").Value = 2999
$ws.Range("K89").Value = 14995
$ws.Range("M89").Value = -9379

$ws.Range("H99").Value = 3001.6
$ws.Range("I99").Value = 2264.6
$ws.Range("J99").Value = 3738.6
$ws.Range("K99").Value = 2264.6
$ws.Range("L99").Value = 3738.6
$ws.Range("M99").Value = -766.5999999999999
$ws.Range("N99").Value = -6734.6

$ws.Range("H105").Value = 100027580
$ws.Range("I105").Value = 125033336
$ws.Range("K105").Value = 125033336
$ws.Range("M105").Value = -125031589

$ws.Range("H107").Value = 1918.4706
$ws.Range("I107").Value = 1918.4706
$ws.Range("K107").Value = 1918.4706
$ws.Range("M107").Value = 1.529399999999896

$ws.Range("H134").Value = 5944.423
$ws.Range("I134").Value = 2494.1333
$ws.Range("J134").Value = 10649.363
$ws.Range("K134").Value = 7482.3999
$ws.Range("L134").Value = 31948.089
$ws.Range("M134").Value = -4947.3999
$ws.Range("N134").Value = -37018.089

$ws.Range("H138").Value = 87711.2
$ws.Range("J138").Value = 87711.2
$ws.Range("L138").Value = 87711.2
$ws.Range("N138").Value = -97991.2


$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 35718210
$ws.Range("I31").Value = 90910000
$ws.Range("K31").Value = 90910000
$ws.Range("M31").Value = -90909705

$ws.Range("H34").Value = 35718210
$ws.Range("I34").Value = 90910000
$ws.Range("K34").Value = 90910000
$ws.Range("M34").Value = -90909798

$ws.Range("H45").Value = 14666.333
$ws.Range("J45").Value = 14999.5
$ws.Range("L45").Value = 14999.5
$ws.Range("N45").Value = -16185.5

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H105").Value = 2139.875
$ws.Range("I105").Value = 2139.875
$ws.Range("K105").Value = 2139.875
$ws.Range("M105").Value = -392.875

$ws.Range("H132").Value = 68356.10000000001
$ws.Range("I132").Value = 6372.8335
$ws.Range("K132").Value = 19118.5005
$ws.Range("M132").Value = -16588.5005


$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 33053854
$ws.Range("J4").Value = 101480.1
$ws.Range("L4").Value = 304440.3
$ws.Range("N4").Value = -304664.3

$ws.Range("H17").Value = 485.14285
$ws.Range("I17").Value = 485.14285
$ws.Range("K17").Value = 1455.42855
$ws.Range("M17").Value = -1286.42855

$ws.Range("H131").Value = 13892707
$ws.Range("I131").Value = 41667596
$ws.Range("J131").Value = 5262.4375
$ws.Range("K131").Value = 125002788
$ws.Range("L131").Value = 15787.3125
$ws.Range("M131").Value = -124997748
$ws.Range("N131").Value = -25867.3125

$ws.Range("H139").Value = 83334970
$ws.Range("I139").Value = 83334970
$ws.Range("K139").Value = 250004910
$ws.Range("M139").Value = -249999770

$ws.Range("H140").Value = 942.7857
$ws.Range("I140").Value = 633.2222
$ws.Range("K140").Value = 1899.6666
$ws.Range("M140").Value = 3280.3334

$ws.Range("H141").Value = 3884.8
$ws.Range("I141").Value = 3884.8
$ws.Range("K141").Value = 11654.4
$ws.Range("M141").Value = -6474.400000000001


$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H14").Value = 809450.7
$ws.Range("J14").Value = 676669.3
$ws.Range("L14").Value = 676669.3
$ws.Range("N14").Value = -677005.3

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H33").Value = 11499.625
$ws.Range("I33").Value = 8000
$ws.Range("J33").Value = 21998.5
$ws.Range("K33").Value = 8000
$ws.Range("L33").Value = 21998.5
$ws.Range("M33").Value = -7748
$ws.Range("N33").Value = -22502.5

$ws.Range("H70").Value = 16893.291
$ws.Range("I70").Value = 14005.046
$ws.Range("K70").Value = 14005.046
$ws.Range("M70").Value = -13735.046

$ws.Range("H73").Value = 16893.291
$ws.Range("I73").Value = 14005.046
$ws.Range("K73").Value = 14005.046
$ws.Range("M73").Value = -13069.046

$ws.Range("H97").Value = 2166.6667
$ws.Range("J97").Value = 1800
$ws.Range("L97").Value = 1800
$ws.Range("N97").Value = -2792

$ws.Range("H107").Value = 736
$ws.Range("I107").Value = 438.8
$ws.Range("K107").Value = 438.8
$ws.Range("M107").Value = 1481.2

$ws.Range("H126").Value = 2235.7778
$ws.Range("I126").Value = 2235.7778
$ws.Range("K126").Value = 6707.3334
$ws.Range("M126").Value = -4237.3334

$ws.Range("H132").Value = 5641.9697
$ws.Range("I132").Value = 4817.5
$ws.Range("J132").Value = 13886.667
$ws.Range("K132").Value = 14452.5
$ws.Range("L132").Value = 41660.001
$ws.Range("M132").Value = -11922.5
$ws.Range("N132").Value = -46720.001


$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H16").Value = 1842.9231
$ws.Range("I16").Value = 1178
$ws.Range("K16").Value = 1178
$ws.Range("M16").Value = -1008

$ws.Range("H122").Value = 5322.905
$ws.Range("I122").Value = 4735.143
$ws.Range("J122").Value = 6498.4287
$ws.Range("K122").Value = 14205.429
$ws.Range("L122").Value = 19495.2861
$ws.Range("M122").Value = -11755.429
$ws.Range("N122").Value = -24395.2861


$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H127").Value = 64998
$ws.Range("J127").Value = 64998
$ws.Range("L127").Value = 64998
$ws.Range("N127").Value = -74918

$ws.Range("H136").Value = 7747.6665
$ws.Range("I136").Value = 3995
$ws.Range("K136").Value = 11985
$ws.Range("M136").Value = -9435

